$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New "Zeitprotokoll" header block (columns E:H) ------------------------

# Column widths for the new block (closest representable value to 20.7109375)
$ws.Range("E1").ColumnWidth = 19.85
$ws.Range("F1").ColumnWidth = 19.85
$ws.Range("G1").ColumnWidth = 19.85
$ws.Range("H1").ColumnWidth = 19.85

# Row 2 / Row 3 "Output"-style label + value cells
$ws.Range("E2:G3").Style = "Output"

$ws.Range("H2").Style = "Output"
$ws.Range("H2").HorizontalAlignment = -4152   # xlRight

$ws.Range("H3").Style = "Calculation"
$style = $ws.Range("H3").Style
$style.NumberFormat = "dd:hh:mm"
$ws.Range("H3").Formula = "=SUM(C:C)"

# Cell text content (order chosen so that shared-string indices line up with
# the target workbook: Name, Vorname, Projekt, C0/H0 Transformation,
# gesamte Arbeitsdauer, Zeitprotokoll)
$ws.Range("E2").Value = "Name:"
$ws.Range("E3").Value = "Vorname:"
$ws.Range("G2").Value = "Projekt:"
$ws.Range("H2").Value = "C0/H0 Transformation"
$ws.Range("G3").Value = "gesamte Arbeitsdauer:"

# Title row, merged across E1:H1
$ws.Range("E1:H1").Merge()
$ws.Range("E1").Value = "Zeitprotokoll"
$ws.Range("E1:H1").Style = "Check Cell"
$ws.Range("E1:H1").HorizontalAlignment = -4108   # xlCenter

# Row 1 is now visually taller (bold title text)
$ws.Rows(1).RowHeight = 16.5

# Selection ends on D6, like in the edited workbook
$ws.Range("D6").Select()
